$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(13).Delete()

$ws.Range("B15").Value = "Semestral-test"
$ws.Range("C15").Value = "Semestral-test"

Write-Host "done"
